# Quarterly indexing esoteric bug-fix operation
# Updates the computed error-statistics table (rows 2-11, columns B:G)
# on the "first_eval" sheet with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> values for columns B,C,D,E,F,G (ME, MAE, MSE, RMSE, SE, N)
$rows = @{
    2  = @(-0.3072063324736731, 0.5773635797369776, 0.6188915608581815, 0.7866966129698166, 0.733013573663008,  42)
    3  = @(0.06855553159162397, 0.8104126323077684, 1.411859963597413,  1.188217136552664,  1.200974231718875,  41)
    4  = @(-0.2062133925061186, 0.5585097485030701, 0.5025698948149058, 0.7089216422249399, 0.6869075491602615, 40)
    5  = @(-0.02557072480903762,0.6736167808661312, 0.868848526268419,  0.9321204462237801, 0.9439501535086928, 39)
    6  = @(-0.1229089188492291, 0.4460013948249449, 0.350767677791658,  0.592256429084276,  0.5871396736870387, 38)
    7  = @(0.05574729004729425, 0.6872371239517284, 0.9195149290446228, 0.9589134106083942, 0.9704962253316772, 37)
    8  = @(-0.1214114959447097, 0.4048482891640412, 0.2561152986148548, 0.5060783522487944, 0.4982679629511933, 36)
    9  = @(0.1522458912649086,  0.6577382853197382, 0.8416789472939559, 0.9174306226053041, 0.9179181168793588, 35)
    10 = @(-0.03691504020070672,0.3541807863119035, 0.2034857836490605, 0.4510939853833794, 0.4563419822598951, 34)
    11 = @(0.2613592320042264,  0.6789426103936371, 1.020165729568819,  1.010032538866357,  0.9907586059223498, 33)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B - ME
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C - MAE
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D - MSE
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E - RMSE
    $ws.Cells.Item($r, 6).Value = $vals[4]   # F - SE
    $ws.Cells.Item($r, 7).Value = $vals[5]   # G - N
}
